$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" right before "2022-Q2" (2nd tab), which
#    pushes 2022-Q2 / 2022-Q1 / 2021-Q4 / 2021-Q3 / 2021-Q2 one slot to the
#    right.
# ---------------------------------------------------------------------------
$sheetQ2 = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($sheetQ2)
$newSheet.Name = "2022-Q3"

# Copy formatting (fonts/borders/alignment) from the "2022-Q2" sheet (now
# pushed to position 3) so the new sheet matches the look of its siblings.
# (Column A on row 1 is always left blank on these sheets, so the header and
# body ranges are copied separately to avoid manufacturing an A1 cell.)
$refSheet = $wb.Worksheets.Item(3)
$refSheet.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)
$refSheet.Range("A2:H6").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# ---------------------------------------------------------------------------
# 2. Populate "2022-Q3" with the fund-holding table.
# ---------------------------------------------------------------------------
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

$fundRows = @(
    @(0,"005457","景顺长城量化小盘股票","6.57","93.58","1.83","0.1202",2),
    @(1,"015496","景顺中证1000指数增强C","1.83","92.63","1.82","0.0333",3),
    @(2,"012879","中信建投量化精选6个月持有期混合C","3.33","90.73","0.99","0.0330",9),
    @(3,"012878","中信建投量化精选6个月持有期混合A","1.67","90.73","0.99","0.0165",9),
    @(4,"015495","景顺中证1000指数增强A","0.69","92.63","1.82","0.0126",3)
)

$textCols = @(2,3,4,5,6,7)
$rowIdx = 2
foreach ($fr in $fundRows) {
    $newSheet.Cells.Item($rowIdx,1).Value = $fr[0]
    foreach ($col in $textCols) {
        $cell = $newSheet.Cells.Item($rowIdx,$col)
        $cell.NumberFormat = "@"
        $cell.Value = $fr[$col - 1]
        $cell.Style = "Normal"
    }
    $newSheet.Cells.Item($rowIdx,8).Value = $fr[7]
    $rowIdx = $rowIdx + 1
}

# ---------------------------------------------------------------------------
# 3. Update the "总计" summary sheet: shift the existing 5 rows down by one
#    and insert the new "2022-Q3" row at the top.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

for ($r = 6; $r -ge 2; $r--) {
    $dst = $r + 1
    $summary.Cells.Item($dst,2).Value = $summary.Cells.Item($r,2).Value2
    $summary.Cells.Item($dst,3).Value = $summary.Cells.Item($r,3).Value2
    $summary.Cells.Item($dst,4).Value = $summary.Cells.Item($r,4).Value2
}

# Row 7 is brand new - copy the index cell's look from A2 and set its index.
$summary.Range("A2").Copy()
$summary.Range("A7").PasteSpecial(-4122)
$summary.Range("A7").Value = 5

$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 5
$summary.Cells.Item(2,4).Value = 0.22

# ---------------------------------------------------------------------------
# 4. Leave selection/active sheet on "总计", matching the original workbook.
# ---------------------------------------------------------------------------
$summary.Activate()
$summary.Range("A1").Select()
